$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("P2").Value = 1.82
$ws.Range("R2").Value = 1.3
$ws.Range("AG2").Value = 10

# Row 3
$ws.Range("AB3").Value = 16

# Row 4
$ws.Range("X4").Value = 9.4
$ws.Range("AB4").Value = 15.5
$ws.Range("AO4").Value = 16.5

# Row 5
$ws.Range("M5").Value = 1.08
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 1.9
$ws.Range("S5").Value = 3.7
$ws.Range("T5").Value = 1.95
$ws.Range("X5").Value = 14

# Row 6
$ws.Range("F6").Value = 2.82
$ws.Range("G6").Value = 2.94
$ws.Range("H6").Value = 2.64
$ws.Range("I6").Value = 2.74
$ws.Range("T6").Value = 1.77

# Row 7
$ws.Range("N7").Value = 3.9
$ws.Range("S7").Value = 3.55
$ws.Range("AG7").Value = 14

# Row 8
$ws.Range("F8").Value = 1.63
$ws.Range("G8").Value = 1.64
$ws.Range("I8").Value = 6.6
$ws.Range("K8").Value = 4.5
$ws.Range("P8").Value = 2.14
$ws.Range("R8").Value = 1.43
$ws.Range("S8").Value = 2.96
$ws.Range("T8").Value = 1.83
$ws.Range("X8").Value = 18.5
$ws.Range("Y8").Value = 24
$ws.Range("Z8").Value = 55
$ws.Range("AA8").Value = 180
$ws.Range("AB8").Value = 9.4
$ws.Range("AD8").Value = 24
$ws.Range("AF8").Value = 9.8
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 14.5
$ws.Range("AK8").Value = 14.5
$ws.Range("AL8").Value = 34
$ws.Range("AM8").Value = 120
$ws.Range("AN8").Value = 9
$ws.Range("AO8").Value = 100
